$wb = $excel.ActiveWorkbook

# --- summary_statistics ---
$ws = $wb.Worksheets.Item("summary_statistics")
$ws.Range("B2").Value = [double]"342"
$ws.Range("C2").Value = [double]"-6.84"
$ws.Range("D2").Value = [double]"7.89"
$ws.Range("E2").Value = [double]"1.19"
$ws.Range("F2").Value = [double]"0.68"
$ws.Range("G2").Value = [double]"1.96"
$ws.Range("H2").Value = [double]"2.27"
$ws.Range("J2").Value = [double]"2.270284413141556"
$ws.Range("B3").Value = [double]"342"
$ws.Range("B4").Value = [double]"342"
$ws.Range("B5").Value = [double]"342"
$ws.Range("B6").Value = [double]"342"
$ws.Range("E6").Value = [double]"0.15"
$ws.Range("B7").Value = [double]"342"
$ws.Range("E7").Value = [double]"0.41"
$ws.Range("B8").Value = [double]"342"
$ws.Range("B9").Value = [double]"342"
$ws.Range("E9").Value = [double]"0.32"
$ws.Range("G9").Value = [double]"0.47"
$ws.Range("B10").Value = [double]"342"
$ws.Range("G10").Value = [double]"0.44"
$ws.Range("H10").Value = [double]"1"
$ws.Range("J10").Value = [double]"1"
$ws.Range("B11").Value = [double]"342"
$ws.Range("G11").Value = [double]"0.21"
$ws.Range("B12").Value = [double]"342"
$ws.Range("E12").Value = [double]"0.33"
$ws.Range("B13").Value = [double]"342"
$ws.Range("E13").Value = [double]"826.27"
$ws.Range("F13").Value = [double]"188.64"
$ws.Range("G13").Value = [double]"2375.13"
$ws.Range("H13").Value = [double]"503.6"
$ws.Range("I13").Value = [double]"25.73816666666667"
$ws.Range("B14").Value = [double]"338"
$ws.Range("E14").Value = [double]"1.57"
$ws.Range("G14").Value = [double]"8.93"
$ws.Range("B15").Value = [double]"338"
$ws.Range("E15").Value = [double]"4.04"
$ws.Range("G15").Value = [double]"7.78"
$ws.Range("B16").Value = [double]"338"
$ws.Range("E16").Value = [double]"13.34"
$ws.Range("G16").Value = [double]"21.61"
$ws.Range("H16").Value = [double]"14.95"
$ws.Range("J16").Value = [double]"15.65"
$ws.Range("B17").Value = [double]"338"
$ws.Range("E17").Value = [double]"25.46"
$ws.Range("F17").Value = [double]"14.65"
$ws.Range("G17").Value = [double]"26.58"
$ws.Range("H17").Value = [double]"30.5"
$ws.Range("J17").Value = [double]"37"
$ws.Range("B18").Value = [double]"338"
$ws.Range("E18").Value = [double]"0.27"
$ws.Range("G18").Value = [double]"0.83"
$ws.Range("B19").Value = [double]"338"
$ws.Range("B20").Value = [double]"338"
$ws.Range("E20").Value = [double]"4.12"
$ws.Range("G20").Value = [double]"3.46"
$ws.Range("H20").Value = [double]"4.28"
$ws.Range("J20").Value = [double]"5.775"
$ws.Range("B21").Value = [double]"338"
$ws.Range("E21").Value = [double]"51.19"
$ws.Range("F21").Value = [double]"55"
$ws.Range("G21").Value = [double]"31.77"
$ws.Range("H21").Value = [double]"57.45"
$ws.Range("I21").Value = [double]"23.15"
$ws.Range("J21").Value = [double]"80.60000000000001"
$ws.Range("B22").Value = [double]"301"
$ws.Range("I22").Value = [double]"0.2271662763466042"
$ws.Range("J22").Value = [double]"0.6909975669099757"
$ws.Range("B23").Value = [double]"342"
$ws.Range("E23").Value = [double]"51.31"
$ws.Range("F23").Value = [double]"52.9"
$ws.Range("G23").Value = [double]"27.83"
$ws.Range("H23").Value = [double]"48.49"
$ws.Range("I23").Value = [double]"26.89"

# --- dichotomous_stats ---
$ws = $wb.Worksheets.Item("dichotomous_stats")
$ws.Range("B2").Value = [double]"302"
$ws.Range("D2").Value = [double]"0.28"
$ws.Range("E2").Value = [double]"1.224"
$ws.Range("G2").Value = [double]"1.05"
$ws.Range("H2").Value = [double]"0.298"
$ws.Range("I2").Value = [double]"58.629"
$ws.Range("J2").Value = [double]"-0.254"
$ws.Range("K2").Value = [double]"0.8139999999999999"
$ws.Range("B3").Value = [double]"253"
$ws.Range("D3").Value = [double]"0.271"
$ws.Range("E3").Value = [double]"1.262"
$ws.Range("G3").Value = [double]"1.091"
$ws.Range("H3").Value = [double]"0.277"
$ws.Range("I3").Value = [double]"146.578"
$ws.Range("J3").Value = [double]"-0.219"
$ws.Range("K3").Value = [double]"0.76"
$ws.Range("B4").Value = [double]"252"
$ws.Range("D4").Value = [double]"0.423"
$ws.Range("E4").Value = [double]"1.303"
$ws.Range("G4").Value = [double]"1.884"
$ws.Range("H4").Value = [double]"0.061"
$ws.Range("I4").Value = [double]"178.01"
$ws.Range("J4").Value = [double]"-0.02"
$ws.Range("K4").Value = [double]"0.866"
$ws.Range("B5").Value = [double]"292"
$ws.Range("D5").Value = [double]"0.31"
$ws.Range("E5").Value = [double]"1.237"
$ws.Range("G5").Value = [double]"0.988"
$ws.Range("H5").Value = [double]"0.327"
$ws.Range("I5").Value = [double]"64.544"
$ws.Range("J5").Value = [double]"-0.317"
$ws.Range("K5").Value = [double]"0.9370000000000001"
$ws.Range("B6").Value = [double]"203"
$ws.Range("C6").Value = [double]"139"
$ws.Range("D6").Value = [double]"0.304"
$ws.Range("E6").Value = [double]"1.315"
$ws.Range("F6").Value = [double]"1.011"
$ws.Range("G6").Value = [double]"1.429"
$ws.Range("H6").Value = [double]"0.154"
$ws.Range("I6").Value = [double]"306.756"
$ws.Range("J6").Value = [double]"-0.115"
$ws.Range("K6").Value = [double]"0.723"
$ws.Range("B7").Value = [double]"238"
$ws.Range("D7").Value = [double]"0.209"
$ws.Range("E7").Value = [double]"1.255"
$ws.Range("G7").Value = [double]"0.955"
$ws.Range("H7").Value = [double]"0.341"
$ws.Range("I7").Value = [double]"220.347"
$ws.Range("J7").Value = [double]"-0.223"
$ws.Range("K7").Value = [double]"0.641"
$ws.Range("B8").Value = [double]"233"
$ws.Range("D8").Value = [double]"0.312"
$ws.Range("E8").Value = [double]"1.291"
$ws.Range("G8").Value = [double]"1.404"
$ws.Range("H8").Value = [double]"0.162"
$ws.Range("I8").Value = [double]"221.665"
$ws.Range("J8").Value = [double]"-0.126"
$ws.Range("K8").Value = [double]"0.75"
$ws.Range("B9").Value = [double]"255"
$ws.Range("D9").Value = [double]"0.014"
$ws.Range("E9").Value = [double]"1.195"
$ws.Range("G9").Value = [double]"0.061"
$ws.Range("H9").Value = [double]"0.952"
$ws.Range("I9").Value = [double]"165.604"
$ws.Range("J9").Value = [double]"-0.44"
$ws.Range("K9").Value = [double]"0.467"
$ws.Range("B10").Value = [double]"327"
$ws.Range("D10").Value = [double]"0.375"
$ws.Range("E10").Value = [double]"1.208"
$ws.Range("G10").Value = [double]"0.831"
$ws.Range("H10").Value = [double]"0.418"
$ws.Range("I10").Value = [double]"15.778"
$ws.Range("J10").Value = [double]"-0.583"
$ws.Range("K10").Value = [double]"1.333"
$ws.Range("B11").Value = [double]"230"
$ws.Range("D11").Value = [double]"0.49"
$ws.Range("E11").Value = [double]"1.352"
$ws.Range("G11").Value = [double]"2.329"
$ws.Range("H11").Value = [double]"0.021"
$ws.Range("I11").Value = [double]"259.574"
$ws.Range("J11").Value = [double]"0.076"
$ws.Range("K11").Value = [double]"0.904"

# --- anovas ---
$ws = $wb.Worksheets.Item("anovas")
$ws.Range("C2").Value = [double]"113.43058926293"
$ws.Range("D2").Value = [double]"37.81019642097666"
$ws.Range("E2").Value = [double]"10.73490303293985"
$ws.Range("F2").Value = [double]"9.33618882495955e-07"
$ws.Range("C3").Value = [double]"6.435566734569218"
$ws.Range("D3").Value = [double]"2.145188911523073"
$ws.Range("E3").Value = [double]"0.5588281554436855"
$ws.Range("F3").Value = [double]"0.6425484902713319"
$ws.Range("C4").Value = [double]"351.5793820450792"
$ws.Range("D4").Value = [double]"18.50417800237259"
$ws.Range("E4").Value = [double]"6.256492055995096"
$ws.Range("F4").Value = [double]"1.110040335505228e-13"

# --- continuous_correlations ---
$ws = $wb.Worksheets.Item("continuous_correlations")
$ws.Range("B2").Value = [double]"0.005"
$ws.Range("C2").Value = [double]"0.1"
$ws.Range("D2").Value = [double]"0.921"
$ws.Range("E2").Value = [double]"340"
$ws.Range("F2").Value = [double]"-0.101"
$ws.Range("G2").Value = [double]"0.111"
$ws.Range("C3").Value = [double]"-0.739"
$ws.Range("D3").Value = [double]"0.461"
$ws.Range("E3").Value = [double]"336"
$ws.Range("F3").Value = [double]"-0.146"
$ws.Range("G3").Value = [double]"0.067"
$ws.Range("B4").Value = [double]"0.012"
$ws.Range("C4").Value = [double]"0.221"
$ws.Range("D4").Value = [double]"0.825"
$ws.Range("E4").Value = [double]"336"
$ws.Range("F4").Value = [double]"-0.095"
$ws.Range("G4").Value = [double]"0.119"
$ws.Range("B5").Value = [double]"0.034"
$ws.Range("C5").Value = [double]"0.618"
$ws.Range("D5").Value = [double]"0.537"
$ws.Range("E5").Value = [double]"336"
$ws.Range("F5").Value = [double]"-0.073"
$ws.Range("G5").Value = [double]"0.14"
$ws.Range("B6").Value = [double]"0.067"
$ws.Range("C6").Value = [double]"1.224"
$ws.Range("D6").Value = [double]"0.222"
$ws.Range("E6").Value = [double]"336"
$ws.Range("F6").Value = [double]"-0.04"
$ws.Range("G6").Value = [double]"0.172"
$ws.Range("B7").Value = [double]"-0.032"
$ws.Range("C7").Value = [double]"-0.596"
$ws.Range("D7").Value = [double]"0.552"
$ws.Range("E7").Value = [double]"336"
$ws.Range("F7").Value = [double]"-0.139"
$ws.Range("G7").Value = [double]"0.074"
$ws.Range("B8").Value = [double]"-0.08"
$ws.Range("C8").Value = [double]"-1.465"
$ws.Range("D8").Value = [double]"0.144"
$ws.Range("E8").Value = [double]"336"
$ws.Range("F8").Value = [double]"-0.185"
$ws.Range("G8").Value = [double]"0.027"
$ws.Range("B9").Value = [double]"0.006"
$ws.Range("C9").Value = [double]"0.102"
$ws.Range("D9").Value = [double]"0.919"
$ws.Range("E9").Value = [double]"336"
$ws.Range("F9").Value = [double]"-0.101"
$ws.Range("G9").Value = [double]"0.112"
$ws.Range("B10").Value = [double]"-0.07000000000000001"
$ws.Range("C10").Value = [double]"-1.28"
$ws.Range("E10").Value = [double]"336"
$ws.Range("F10").Value = [double]"-0.175"
$ws.Range("B11").Value = [double]"0.016"
$ws.Range("C11").Value = [double]"0.282"
$ws.Range("D11").Value = [double]"0.778"
$ws.Range("E11").Value = [double]"299"
$ws.Range("F11").Value = [double]"-0.097"
$ws.Range("G11").Value = [double]"0.129"
$ws.Range("B12").Value = [double]"0.062"
$ws.Range("C12").Value = [double]"1.147"
$ws.Range("D12").Value = [double]"0.252"
$ws.Range("E12").Value = [double]"340"
$ws.Range("F12").Value = [double]"-0.044"
$ws.Range("G12").Value = [double]"0.167"
